# Build site at 2022-09-26 16:07:08 UTC
#
# The LOQ4212 syllabus sheet had its stray "Docentes responsaveis:" value
# row (old row 13, holding only "5840917 - Fabricio Maciel Gomes" with no
# label of its own) removed, which shifts every row below it up by one
# (new dimension A1:C23). On top of that shift, a handful of B/C value
# cells were overwritten with text already used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old stand-alone "5840917 - Fabricio Maciel Gomes" row; rows
# 14-24 shift up to become rows 13-23.
$ws.Rows.Item(13).Delete()

# Patch the B/C (value) columns that no longer line up with their
# previous contents after the shift.
$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2018" already exists elsewhere on the sheet as text (B8/C8) -
# copy it across instead of assigning the literal string, which Excel
# would otherwise auto-convert into a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.Application.CutCopyMode = $false

$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios."

$ws.Range("B20").Value = "MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas."
$ws.Range("C20").Value = "MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas."

$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

$wb.Save()
